$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the A2:B5 values with the new cluster data
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1496

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 1142

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 949

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 658

# Remove the now-obsolete row 6 entirely
$ws.Rows.Item(6).Delete()
